$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("F16").Value = "Al Nasr"
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = "Al Sharjah"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 3.33
$ws.Range("K16").Value = "16/09/2023 18:13"
$ws.Range("L16").Value = 3.73
$ws.Range("M16").Value = "23/09/2023 15:11"
$ws.Range("N16").Value = 3.6
$ws.Range("O16").Value = "16/09/2023 18:13"
$ws.Range("P16").Value = 3.93
$ws.Range("Q16").Value = "23/09/2023 15:17"
$ws.Range("R16").Value = 2.1
$ws.Range("S16").Value = "16/09/2023 18:13"
$ws.Range("T16").Value = 1.91
$ws.Range("U16").Value = "23/09/2023 15:17"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-nasr-al-sharjah/bytlvyiD/"

# Row 17
$ws.Range("F17").Value = "Ittihad Kalba"
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = "Al Jazira"
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 3.33
$ws.Range("K17").Value = "16/09/2023 18:13"
$ws.Range("L17").Value = 4.23
$ws.Range("M17").Value = "23/09/2023 15:11"
$ws.Range("N17").Value = 3.77
$ws.Range("O17").Value = "16/09/2023 18:13"
$ws.Range("P17").Value = 4.42
$ws.Range("Q17").Value = "23/09/2023 15:11"
$ws.Range("R17").Value = 2.04
$ws.Range("S17").Value = "16/09/2023 18:13"
$ws.Range("T17").Value = 1.71
$ws.Range("U17").Value = "23/09/2023 15:11"
$ws.Range("V17").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/ittihad-kalba-al-jazira/x6s0ohbt/"

# Row 20
$ws.Range("F20").Value = "Al Wahda"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "Hatta"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1.36
$ws.Range("K20").Value = "19/09/2023 16:42"
$ws.Range("L20").Value = 1.24
$ws.Range("M20").Value = "24/09/2023 14:22"
$ws.Range("N20").Value = 5.48
$ws.Range("O20").Value = "19/09/2023 16:42"
$ws.Range("P20").Value = 6.54
$ws.Range("Q20").Value = "24/09/2023 15:15"
$ws.Range("R20").Value = 6.99
$ws.Range("S20").Value = "19/09/2023 16:42"
$ws.Range("T20").Value = 10.02
$ws.Range("U20").Value = "24/09/2023 15:15"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wahda-hatta/pWp8qWSh/"

# Row 21
$ws.Range("F21").Value = "Al Bataeh"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "Shabab Al-Ahli Dubai"
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 6.4
$ws.Range("K21").Value = "17/09/2023 15:42"
$ws.Range("L21").Value = 6.67
$ws.Range("M21").Value = "24/09/2023 15:15"
$ws.Range("N21").Value = 5.03
$ws.Range("O21").Value = "17/09/2023 15:42"
$ws.Range("P21").Value = 4.85
$ws.Range("Q21").Value = "24/09/2023 15:15"
$ws.Range("R21").Value = 1.37
$ws.Range("S21").Value = "17/09/2023 15:42"
$ws.Range("T21").Value = 1.44
$ws.Range("U21").Value = "24/09/2023 15:15"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-shabab-al-ahli-dubai/6un4pCDn/"

# Row 52
$ws.Range("F52").Value = "Hatta"
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = "Bani Yas"
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 3.08
$ws.Range("K52").Value = "24/11/2023 15:12"
$ws.Range("L52").Value = 4.35
$ws.Range("M52").Value = "25/11/2023 13:43"
$ws.Range("N52").Value = 3.72
$ws.Range("O52").Value = "24/11/2023 15:12"
$ws.Range("P52").Value = 4.42
$ws.Range("Q52").Value = "25/11/2023 13:43"
$ws.Range("R52").Value = 2.09
$ws.Range("S52").Value = "24/11/2023 15:12"
$ws.Range("T52").Value = 1.69
$ws.Range("U52").Value = "25/11/2023 13:43"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/hatta-bani-yas/6uE7LPuA/"

# Row 53
$ws.Range("F53").Value = "Al Wahda"
$ws.Range("G53").Value = 4
$ws.Range("H53").Value = "Emirates Club"
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 1.33
$ws.Range("K53").Value = "24/11/2023 15:12"
$ws.Range("L53").Value = 1.31
$ws.Range("M53").Value = "25/11/2023 13:36"
$ws.Range("N53").Value = 5.27
$ws.Range("O53").Value = "24/11/2023 15:12"
$ws.Range("P53").Value = 5.9
$ws.Range("Q53").Value = "25/11/2023 13:36"
$ws.Range("R53").Value = 8.22
$ws.Range("S53").Value = "24/11/2023 15:12"
$ws.Range("T53").Value = 8.33
$ws.Range("U53").Value = "25/11/2023 13:36"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wahda-emirates-club/63dfCL2d/"

# Row 67
$ws.Range("F67").Value = "Hatta"
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = "Al Sharjah"
$ws.Range("I67").Value = 4
$ws.Range("J67").Value = 18.05
$ws.Range("K67").Value = "08/12/2023 14:35"
$ws.Range("L67").Value = 10.68
$ws.Range("M67").Value = "09/12/2023 13:43"
$ws.Range("N67").Value = 9.18
$ws.Range("O67").Value = "08/12/2023 14:35"
$ws.Range("P67").Value = 6.38
$ws.Range("Q67").Value = "09/12/2023 13:43"
$ws.Range("R67").Value = 1.1
$ws.Range("S67").Value = "08/12/2023 14:35"
$ws.Range("T67").Value = 1.24
$ws.Range("U67").Value = "09/12/2023 13:43"
$ws.Range("V67").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/hatta-al-sharjah/j9wUYZFQ/"

# Row 68
$ws.Range("F68").Value = "Khorfakkan"
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = "Al Nasr"
$ws.Range("I68").Value = 1
$ws.Range("J68").Value = 3.39
$ws.Range("K68").Value = "05/12/2023 15:42"
$ws.Range("L68").Value = 5.71
$ws.Range("M68").Value = "09/12/2023 13:43"
$ws.Range("N68").Value = 3.69
$ws.Range("O68").Value = "05/12/2023 15:42"
$ws.Range("P68").Value = 4.73
$ws.Range("Q68").Value = "09/12/2023 13:43"
$ws.Range("R68").Value = 1.97
$ws.Range("S68").Value = "05/12/2023 15:42"
$ws.Range("T68").Value = 1.51
$ws.Range("U68").Value = "09/12/2023 13:43"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/khorfakkan-al-nasr/UgkLzdoE/"

# New row 78 - copy formatting from row 77 first, then set values
$ws.Range("A77").Copy($ws.Range("A78"))
$ws.Range("E77").Copy($ws.Range("E78"))

$ws.Range("A78").Value = 77
$ws.Range("B78").Value = "united-arab-emirates"
$ws.Range("C78").Value = "uae-league"
$ws.Range("D78").Value = "2023-2024"
$ws.Range("E78").Value = 45280.6875
$ws.Range("F78").Value = "Al Sharjah"
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = "Shabab Al-Ahli Dubai"
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = 2.51
$ws.Range("K78").Value = "15/12/2023 16:42"
$ws.Range("L78").Value = 2.43
$ws.Range("M78").Value = "20/12/2023 16:29"
$ws.Range("N78").Value = 3.45
$ws.Range("O78").Value = "15/12/2023 16:42"
$ws.Range("P78").Value = 3.77
$ws.Range("Q78").Value = "20/12/2023 16:29"
$ws.Range("R78").Value = 2.59
$ws.Range("S78").Value = "15/12/2023 16:42"
$ws.Range("T78").Value = 2.72
$ws.Range("U78").Value = "20/12/2023 16:27"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-sharjah-shabab-al-ahli-dubai/88zwy5uc/"
